$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append " 今天心情就这样" to the "2023.3.22哈哈哈哈" paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute("2023.3.22哈哈哈哈", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023.3.22哈哈哈哈 今天心情就这样", 2)

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of that paragraph to the
#    end of the very first paragraph ("2023年3月1日").
# ------------------------------------------------------------------

# Remove the existing _GoBack bookmark (currently sitting after the
# "2023.3.22..." run).
$d.Bookmarks("_GoBack").Delete()

# Re-create it right after the run in paragraph 1, without splitting
# that run. We do this by temporarily inserting a placeholder
# character right after the run, wrapping *that* (non-zero-width)
# range with the bookmark, and then deleting the placeholder again -
# the bookmark collapses cleanly to sit between the run and the
# paragraph mark.
$p1 = $d.Paragraphs(1).Range
$insertPos = $p1.End - 1
$tmp = $d.Range($insertPos, $insertPos)
$tmp.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $tmp)
$placeholder = $d.Range($insertPos, $insertPos + 1)
$placeholder.Text = ""

# ------------------------------------------------------------------
# 3. Styles: mark the "Normal Table" style as a Quick Style
#    (<w:qFormat/> on the style definition).
# ------------------------------------------------------------------
$tableStyle = $d.Styles.Item("Normal Table")
$tableStyle.QuickStyle = $true
